$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")

# New shared strings used by column O (the "campo / control name" column)
# and the new Q-column concat formulas for the Produto table (rows 21-31).

# Row 21: id_produto -> txtCodigo (already-existing shared string)
$ws.Range("O21").Value = "txtCodigo"
$ws.Range("Q21").Formula = '=""""&B21&"_"&$C$19&","""&" +"'

# Row 22: id_categoria_produto -> txtCodigoCategoria (new string)
$ws.Range("O22").Value = "txtCodigoCategoria"

# Row 23: nome_produto -> txtNome (already-existing shared string)
$ws.Range("O23").Value = "txtNome"

# Row 24: qtde_produto -> txtQtde (new string)
$ws.Range("O24").Value = "txtQtde"

# Row 25: peso_produto -> txtPeso (new string)
$ws.Range("O25").Value = "txtPeso"

# Row 26: unidade_produto -> txtUnidade (new string)
$ws.Range("O26").Value = "txtUnidade"

# Row 27: cadastro_produto -> txtCadastro (new string)
$ws.Range("O27").Value = "txtCadastro"

# Row 28: valorCusto_produto -> txtCusto (new string)
$ws.Range("O28").Value = "txtCusto"

# Row 29: valorVenda_produto -> txtVenda (new string)
$ws.Range("O29").Value = "txtVenda"

# Row 30: status_produto -> cboStatus (already-existing shared string)
$ws.Range("O30").Value = "cboStatus"

# Row 31: obs_produto -> txtObs (already-existing shared string)
$ws.Range("O31").Value = "txtObs"

# Rows 22-31 share the same relative concat formula; fill the whole block at
# once so Excel records it as a single shared-formula group (matches how the
# Q-column ranges above it - Q4:Q9, Q14:Q17 - were originally authored).
$ws.Range("Q22:Q31").Formula = '=""""&B22&"_"&$C$19&","""&" +"'

# Reflect the final UI state left behind after the edit: the user scrolled
# so row 10 is at the top and ended with Q21 selected.
$ws.Range("Q21").Select()
$excel.ActiveWindow.ScrollRow = 10
